# Update scripts with new TPM (transcripts-per-million) values for the
# Dll1-Notch1 ligand-receptor pair sheet, and drop the "Resolving-Mac"
# sending-cluster rows that no longer apply with the refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last 4 rows (sending cluster = Resolving-Mac) are removed entirely;
# deleting the range shifts the remaining rows up and shrinks the used range
# from A1:T17 down to A1:T13.
$ws.Range("A14:T17").Delete()

# Row 2: Sending cluster = ECs, Target cluster = ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dll1"
$ws.Range("C2").Value = "Notch1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 16.58130233333334
$ws.Range("H2").Value = 49.74390700000001
$ws.Range("I2").Value = 0.6205214785234227
$ws.Range("J2").Value = 0.6205214785234225
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 58.95713633333333
$ws.Range("N2").Value = 176.871409
$ws.Range("O2").Value = 0.4863146960083892
$ws.Range("P2").Value = 0.4863146960083893
$ws.Range("Q2").Value = 977.5861022505516
$ws.Range("R2").Value = 8798.274920254964
$ws.Range("S2").Value = 0.3017687141947945
$ws.Range("T2").Value = 0.3017687141947945

# Row 3: Sending cluster = ECs, Target cluster = FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dll1"
$ws.Range("C3").Value = "Notch1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 16.58130233333334
$ws.Range("H3").Value = 49.74390700000001
$ws.Range("I3").Value = 0.6205214785234227
$ws.Range("J3").Value = 0.6205214785234225
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.99153
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.07416766570679004
$ws.Range("P3").Value = 0.07416766570679005
$ws.Range("Q3").Value = 149.0912773692367
$ws.Range("R3").Value = 1341.82149632313
$ws.Range("S3").Value = 0.0460226295830083
$ws.Range("T3").Value = 0.0460226295830083

# Row 4: Sending cluster = ECs, Target cluster = MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dll1"
$ws.Range("C4").Value = "Notch1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 16.58130233333334
$ws.Range("H4").Value = 49.74390700000001
$ws.Range("I4").Value = 0.6205214785234227
$ws.Range("J4").Value = 0.6205214785234225
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 42.51661933333333
$ws.Range("N4").Value = 127.549858
$ws.Range("O4").Value = 0.3507032073181665
$ws.Range("P4").Value = 0.3507032073181665
$ws.Range("Q4").Value = 704.9809193572452
$ws.Range("R4").Value = 6344.828274215207
$ws.Range("S4").Value = 0.2176188727279751
$ws.Range("T4").Value = 0.2176188727279751

# Row 5: Sending cluster = ECs, Target cluster = Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Dll1"
$ws.Range("C5").Value = "Notch1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 16.58130233333334
$ws.Range("H5").Value = 49.74390700000001
$ws.Range("I5").Value = 0.6205214785234227
$ws.Range("J5").Value = 0.6205214785234225
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.76719366666667
$ws.Range("N5").Value = 32.301581
$ws.Range("O5").Value = 0.0888144309666542
$ws.Range("P5").Value = 0.08881443096665421
$ws.Range("Q5").Value = 178.5340934685519
$ws.Range("R5").Value = 1606.806841216967
$ws.Range("S5").Value = 0.05511126201764472
$ws.Range("T5").Value = 0.05511126201764472

# Row 6: Sending cluster = FAPs, Target cluster = ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dll1"
$ws.Range("C6").Value = "Notch1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5586273333333333
$ws.Range("H6").Value = 1.675882
$ws.Range("I6").Value = 0.02090549052511678
$ws.Range("J6").Value = 0.02090549052511678
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 58.95713633333333
$ws.Range("N6").Value = 176.871409
$ws.Range("O6").Value = 0.4863146960083892
$ws.Range("P6").Value = 0.4863146960083893
$ws.Range("Q6").Value = 32.93506785085977
$ws.Range("R6").Value = 296.415610657738
$ws.Range("S6").Value = 0.01016664726962843
$ws.Range("T6").Value = 0.01016664726962843

# Row 7: Sending cluster = FAPs, Target cluster = FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dll1"
$ws.Range("C7").Value = "Notch1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5586273333333333
$ws.Range("H7").Value = 1.675882
$ws.Range("I7").Value = 0.02090549052511678
$ws.Range("J7").Value = 0.02090549052511678
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.99153
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.07416766570679004
$ws.Range("P7").Value = 0.07416766570679005
$ws.Range("Q7").Value = 5.022914426486666
$ws.Range("R7").Value = 45.20622983838
$ws.Range("S7").Value = 0.001550511432703328
$ws.Range("T7").Value = 0.001550511432703328

# Row 8: Sending cluster = FAPs, Target cluster = MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Dll1"
$ws.Range("C8").Value = "Notch1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5586273333333333
$ws.Range("H8").Value = 1.675882
$ws.Range("I8").Value = 0.02090549052511678
$ws.Range("J8").Value = 0.02090549052511678
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 42.51661933333333
$ws.Range("N8").Value = 127.549858
$ws.Range("O8").Value = 0.3507032073181665
$ws.Range("P8").Value = 0.3507032073181665
$ws.Range("Q8").Value = 23.75094568052844
$ws.Range("R8").Value = 213.758511124756
$ws.Range("S8").Value = 0.007331622577717994
$ws.Range("T8").Value = 0.007331622577717995

# Row 9: Sending cluster = FAPs, Target cluster = Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Dll1"
$ws.Range("C9").Value = "Notch1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5586273333333333
$ws.Range("H9").Value = 1.675882
$ws.Range("I9").Value = 0.02090549052511678
$ws.Range("J9").Value = 0.02090549052511678
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.76719366666667
$ws.Range("N9").Value = 32.301581
$ws.Range("O9").Value = 0.0888144309666542
$ws.Range("P9").Value = 0.08881443096665421
$ws.Range("Q9").Value = 6.014848685493554
$ws.Range("R9").Value = 54.13363816944199
$ws.Range("S9").Value = 0.001856709245067028
$ws.Range("T9").Value = 0.001856709245067028

# Row 10: Sending cluster = MuSCs, Target cluster = ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Dll1"
$ws.Range("C10").Value = "Notch1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.581631
$ws.Range("H10").Value = 28.744893
$ws.Range("I10").Value = 0.3585730309514606
$ws.Range("J10").Value = 0.3585730309514606
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 58.95713633333333
$ws.Range("N10").Value = 176.871409
$ws.Range("O10").Value = 0.4863146960083892
$ws.Range("P10").Value = 0.4863146960083893
$ws.Range("Q10").Value = 564.9055251626929
$ws.Range("R10").Value = 5084.149726464238
$ws.Range("S10").Value = 0.1743793345439663
$ws.Range("T10").Value = 0.1743793345439663

# Row 11: Sending cluster = MuSCs, Target cluster = FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Dll1"
$ws.Range("C11").Value = "Notch1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9.581631
$ws.Range("H11").Value = 28.744893
$ws.Range("I11").Value = 0.3585730309514606
$ws.Range("J11").Value = 0.3585730309514606
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.99153
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.07416766570679004
$ws.Range("P11").Value = 0.07416766570679005
$ws.Range("Q11").Value = 86.15352258543
$ws.Range("R11").Value = 775.38170326887
$ws.Range("S11").Value = 0.02659452469107841
$ws.Range("T11").Value = 0.02659452469107841

# Row 12: Sending cluster = MuSCs, Target cluster = MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Dll1"
$ws.Range("C12").Value = "Notch1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.581631
$ws.Range("H12").Value = 28.744893
$ws.Range("I12").Value = 0.3585730309514606
$ws.Range("J12").Value = 0.3585730309514606
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 42.51661933333333
$ws.Range("N12").Value = 127.549858
$ws.Range("O12").Value = 0.3507032073181665
$ws.Range("P12").Value = 0.3507032073181665
$ws.Range("Q12").Value = 407.378557819466
$ws.Range("R12").Value = 3666.407020375194
$ws.Range("S12").Value = 0.1257527120124734
$ws.Range("T12").Value = 0.1257527120124734

# Row 13: Sending cluster = MuSCs, Target cluster = Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Dll1"
$ws.Range("C13").Value = "Notch1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.581631
$ws.Range("H13").Value = 28.744893
$ws.Range("I13").Value = 0.3585730309514606
$ws.Range("J13").Value = 0.3585730309514606
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.76719366666667
$ws.Range("N13").Value = 32.301581
$ws.Range("O13").Value = 0.0888144309666542
$ws.Range("P13").Value = 0.08881443096665421
$ws.Range("Q13").Value = 103.167276619537
$ws.Range("R13").Value = 928.505489575833
$ws.Range("S13").Value = 0.03184645970394245
$ws.Range("T13").Value = 0.03184645970394246
